$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new row 64: date, hours, and description
$ws.Range("A64").Value = (Get-Date -Year 2023 -Month 11 -Day 28).Date
$ws.Range("A64").NumberFormat = $ws.Range("A63").NumberFormat
$ws.Range("B64").Value = 4
$ws.Range("C64").Value = "intergrating the travel advisory with the front-end was still giving me a lot of errors, so now im working on the integration woth the weather api, but this is giving a lot of errors as well"

# Update the view/selection to match after-edit state
$ws.Application.ActiveWindow.ScrollRow = 53
$ws.Range("C64").Select()
